$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Conjunto de Treinamento")
$ws.Range("B1").EntireColumn.Insert()
Write-Host ($ws.Range("A1").Value2)
Write-Host ($ws.Range("B1").Value2)
Write-Host ($ws.Range("C1").Value2)
